# New groups and video link
$wb = $excel.ActiveWorkbook

# --- Group_Members sheet: add new Group 4 members ---
$wsMembers = $wb.Worksheets.Item("Group_Members")
$wsMembers.Range("A11").Value = 4
$wsMembers.Range("B11").Value = "Austin Nguyen"
$wsMembers.Range("A12").Value = 4
$wsMembers.Range("B12").Value = "Colleen Xu"
$wsMembers.Range("A13").Value = 4
$wsMembers.Range("B13").Value = "Xiao Wang"

# --- Choices sheet: add new Group 4 preferences ---
$wsChoices = $wb.Worksheets.Item("Choices")
$wsChoices.Range("A11").Value = 4
$wsChoices.Range("B11").Value = 1
$wsChoices.Range("C11").Value = 6
$wsChoices.Range("A12").Value = 4
$wsChoices.Range("B12").Value = 2
$wsChoices.Range("C12").Value = 4
$wsChoices.Range("A13").Value = 4
$wsChoices.Range("B13").Value = 3
$wsChoices.Range("C13").Value = 15

# --- Update selections on each sheet ---
$wsMembers.Range("B14").Select()
$wsChoices.Range("C14").Select()

# --- Make Choices the active sheet (was Group_Members) ---
$wsChoices.Activate()
